$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above the current row 1; old rows 1-24 shift to 3-26.
$ws.Range("A1:O2").EntireRow.Insert()

# New row 1 will hold numeric column indices (0-14). Give it the bold /
# thin-bordered / centered header look that used to live on the original
# header row (which is now row 3), by copying that row's formatting over.
$ws.Range("A3:O3").Copy() | Out-Null
$ws.Range("A1:O1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Row 3 (the former header row) keeps its header text, but loses that
# bold/border/center formatting, reverting to the plain default style.
$ws.Range("A3:O3").ClearFormats()

for ($c = 1; $c -le 15; $c++) {
    $ws.Cells.Item(1, $c).Value = $c - 1
}

# New row 2 is blank except for E2 = "Flat Washer".
for ($c = 1; $c -le 15; $c++) {
    $ws.Cells.Item(2, $c).Value = ""
}
$ws.Cells.Item(2, 5).Value = "Flat Washer"
